$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 456; existing rows 456-516 shift down to 457-517.
$ws.Rows.Item(456).Insert()

# Populate the newly inserted row 456 with its data (same as the row that
# used to be at 456, i.e. now at 457, except for D/J/N which differ).
$ws.Range("A456").Value = 5
$ws.Range("B456").Value = "Macroferia Regional de Talca"
$ws.Range("C456").Value = "Maule"
$ws.Range("D456").Value = 45127
$ws.Range("E456").Value = 7
$ws.Range("F456").Value = 100112003
$ws.Range("G456").Value = "Ajo"
$ws.Range("H456").Value = "Chino"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 200
$ws.Range("K456").Value = 20000
$ws.Range("L456").Value = 20000
$ws.Range("M456").Value = 20000
$ws.Range("N456").Value = "`$/malla 10 kilos"
$ws.Range("O456").Value = "China"
$ws.Range("P456").Value = 2000
$ws.Range("Q456").Value = 10
$ws.Range("R456").Value = "Hortaliza"
